# Applies the 30/12/2025 18:38 scrape update (commit: "🚌 141: 30/12 21:38 LP1912+6203+6173")
# to the three schedule-log sheets: LP1912, LP1912-215, 6203-6173.

$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# Sheet "LP1912": columns A(-), B=Hora_Scrap, C=Hora_Llegada,
# D=Linea, E=Minutos, F=Parada, G=Fecha. 22 new rows (498-519).
# ------------------------------------------------------------------
$wsLP1912 = $wb.Worksheets.Item("LP1912")
$wsLP1912.Range("A2").Value = "Última actualización: 30/12/2025 18:38:18"
$wsLP1912.Range("A3").Value = "Total filas: 518"

$lp1912NewRows = @(
    @(498, "18:38:08", "18:40", "15_ABASTO", 2, "LP1912", "30/12/2025"),
    @(499, "18:38:08", "18:41", "16_SANTA ANA", 3, "LP1912", "30/12/2025"),
    @(500, "18:38:08", "18:48", "14X44_ABASTO", 10, "LP1912", "30/12/2025"),
    @(501, "18:38:08", "18:51", "16_SANTA ANA", 13, "LP1912", "30/12/2025"),
    @(502, "18:38:08", "18:52", "15_ABASTO", 14, "LP1912", "30/12/2025"),
    @(503, "18:38:08", "18:56", "10_OLMOS", 18, "LP1912", "30/12/2025"),
    @(504, "18:38:08", "19:01", "16_SANTA ANA", 23, "LP1912", "30/12/2025"),
    @(505, "18:38:08", "19:05", "11_ETCHEVERRY", 27, "LP1912", "30/12/2025"),
    @(506, "18:38:08", "19:08", "23_HERNANDEZ", 30, "LP1912", "30/12/2025"),
    @(507, "18:38:08", "19:12", "10_OLMOS", 34, "LP1912", "30/12/2025"),
    @(508, "18:38:08", "19:21", "26_HERNANDEZ", 43, "LP1912", "30/12/2025"),
    @(509, "18:38:08", "19:28", "15_ABASTO", 50, "LP1912", "30/12/2025"),
    @(510, "18:38:08", "19:40", "215C_EL PATO", 62, "LP1912", "30/12/2025"),
    @(511, "18:38:08", "19:41", "14_ABASTO", 63, "LP1912", "30/12/2025"),
    @(512, "18:38:08", "19:50", "11X44_ETCHEVERRY", 72, "LP1912", "30/12/2025"),
    @(513, "18:38:08", "19:51", "16_P MOR-SANTA ANA", 73, "LP1912", "30/12/2025"),
    @(514, "18:38:08", "19:51", "81_EL PELIGRO", 73, "LP1912", "30/12/2025"),
    @(515, "18:38:08", "19:59", "17_ROMERO", 81, "LP1912", "30/12/2025"),
    @(516, "18:38:08", "20:01", "14_ABASTO", 83, "LP1912", "30/12/2025"),
    @(517, "18:38:08", "20:08", "10_OLMOS", 90, "LP1912", "30/12/2025"),
    @(518, "18:38:08", "20:11", "16_P MOR-167 Y 521", 93, "LP1912", "30/12/2025"),
    @(519, "18:38:08", "20:13", "23_HERNANDEZ", 95, "LP1912", "30/12/2025")
)

foreach ($row in $lp1912NewRows) {
    $r = $row[0]
    $wsLP1912.Range("B$r").Value = $row[1]
    $wsLP1912.Range("C$r").Value = $row[2]
    $wsLP1912.Range("D$r").Value = $row[3]
    $wsLP1912.Range("E$r").Value = $row[4]
    $wsLP1912.Range("F$r").Value = $row[5]
    $wsLP1912.Range("G$r").Value = $row[6]
}

# ------------------------------------------------------------------
# Sheet "LP1912-215": columns A(-), B=Fecha, C=Hora_Scrap,
# D=Hora_Llegada, E=Linea, F=Minutos, G=Parada. 1 new row (34).
# ------------------------------------------------------------------
$wsLP1912_215 = $wb.Worksheets.Item("LP1912-215")
$wsLP1912_215.Range("A2").Value = "Última actualización: 30/12/2025 18:38:18"
$wsLP1912_215.Range("A3").Value = "Total filas: 33"

$wsLP1912_215.Range("B34").Value = "30/12/2025"
$wsLP1912_215.Range("C34").Value = "18:38:08"
$wsLP1912_215.Range("D34").Value = "19:40"
$wsLP1912_215.Range("E34").Value = "215C_EL PATO"
$wsLP1912_215.Range("F34").Value = 62
$wsLP1912_215.Range("G34").Value = "LP1912"

# ------------------------------------------------------------------
# Sheet "6203-6173": columns A(-), B=Fecha, C=Hora_Scrap,
# D=Hora_Llegada, E=Linea, F=Minutos, G=Parada. 3 new rows (67-69).
# ------------------------------------------------------------------
$ws6203 = $wb.Worksheets.Item("6203-6173")
$ws6203.Range("A2").Value = "Última actualización: 30/12/2025 18:38:18"
$ws6203.Range("A3").Value = "Total filas: 68"

$sheet3NewRows = @(
    @(67, "30/12/2025", "18:38:18", "18:52", "215A_LA PLATA", 14, "L6173"),
    @(68, "30/12/2025", "18:38:18", "19:04", "215B_LP-P MOR-1 Y 57", 26, "L6173"),
    @(69, "30/12/2025", "18:38:13", "19:54", "215C_LA PLATA", 76, "L6203")
)

foreach ($row in $sheet3NewRows) {
    $r = $row[0]
    $ws6203.Range("B$r").Value = $row[1]
    $ws6203.Range("C$r").Value = $row[2]
    $ws6203.Range("D$r").Value = $row[3]
    $ws6203.Range("E$r").Value = $row[4]
    $ws6203.Range("F$r").Value = $row[5]
    $ws6203.Range("G$r").Value = $row[6]
}

Write-Host "Applied 30/12/2025 18:38 update: +22 LP1912, +1 LP1912-215, +3 6203-6173 rows"
